# Horarios actualizados Línea 141 - 362
# Updates the three schedule sheets (LP1912, LP1912-215, 6203-6173) with a
# new scrape timestamp (02:56:21 -> 03:24:15), refreshed arrival estimates,
# and two newly-scraped rows on the main LP1912 sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: LP1912
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 03:24:15"
$ws1.Range("A3").Value = "Total filas: 6"

# Row 6: 14_ABASTO
$ws1.Range("A6").Value = "03:24:15"
$ws1.Range("B6").Value = "03:48"
$ws1.Range("D6").Value = 24

# Row 7: 81_EL PELIGRO
$ws1.Range("A7").Value = "03:24:15"
$ws1.Range("D7").Value = 37

# Row 8: 215_ALUAR
$ws1.Range("A8").Value = "03:24:15"
$ws1.Range("B8").Value = "04:32"
$ws1.Range("D8").Value = 68

# Row 9: 11_ETCHEVERRY
$ws1.Range("A9").Value = "03:24:15"
$ws1.Range("D9").Value = 89

# Row 10 (new): 17_ROMERO
$ws1.Range("A10").Value = "03:24:15"
$ws1.Range("B10").Value = "05:16"
$ws1.Range("C10").Value = "17_ROMERO"
$ws1.Range("D10").Value = 112
$ws1.Range("E10").Value = "LP1912"

# Row 11 (new): 23_HERNANDEZ
$ws1.Range("A11").Value = "03:24:15"
$ws1.Range("B11").Value = "05:22"
$ws1.Range("C11").Value = "23_HERNANDEZ"
$ws1.Range("D11").Value = 118
$ws1.Range("E11").Value = "LP1912"

# ---------------------------------------------------------------------
# Sheet 2: LP1912-215
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 03:24:15"

# Row 6: 215_ALUAR
$ws2.Range("A6").Value = "03:24:15"
$ws2.Range("B6").Value = "04:32"
$ws2.Range("D6").Value = 68

# ---------------------------------------------------------------------
# Sheet 3: 6203-6173
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 03:24:15"
